$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates exactly as described by the diff.
# D-column (Price) values are numeric-looking strings that must remain
# stored as text (matching the original inlineStr "Price" column), so we
# force a text number format before assigning them.

$priceCells = [ordered]@{
    "D2" = "248.16"
    "D4" = "5.546"
    "D5" = "0.05637"
    "D6" = "3.395"
    "D7" = "6.474"
    "D8" = "1.074"
    "D9" = "0.8020"
    "D10" = "0.1429"
    "D11" = "0.07322"
    "D12" = "0.03199"
    "D13" = "0.02989"
    "D15" = "0.001663"
    "D16" = "2.973"
    "D17" = "0.04683"
    "D18" = "0.0005915"
    "D19" = "0.006272"
    "D20" = "0.001058"
    "D21" = "0.003833"
    "D22" = "0.0001501"
    "D23" = "0.0004003"
    "D24" = "3.983"
    "D25" = "2.113"
    "D26" = "0.3290"
    "D27" = "0.1292"
    "D40" = "0.04212"
    "D41" = "0.1048"
    "D42" = "0.002972"
    "D43" = "0.003248"
    "D44" = "0.008747"
    "D45" = "0.00005639"
    "D47" = "0.6805"
    "D48" = "0.02714"
}

foreach ($ref in $priceCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $priceCells[$ref]
}

$textCells = [ordered]@{
    "B6" = "GateToken"
    "C6" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "E6" = "5GateTokenGT"
    "B7" = "KuCoinToken"
    "C7" = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
    "E7" = "6KuCoinTokenKCS"
    "B8" = "FTXToken"
    "C8" = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
    "E8" = "7FTXTokenFTT"
    "B9" = "MXToken"
    "C9" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "E9" = "8MXTokenMX"
    "B16" = "MCDex"
    "C16" = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
    "E16" = "15MCDexMCB"
    "B17" = "CoinExToken"
    "C17" = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
    "E17" = "16CoinExTokenCET"
    "B18" = "One"
    "C18" = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
    "E18" = "17OneONE"
    "B19" = "TigerCash"
    "C19" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "E19" = "18TigerCashTCH"
    "B20" = "BitKan"
    "C20" = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
    "E20" = "19BitKanKAN"
    "B21" = "HotbitToken"
    "C21" = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
    "E21" = "20HotbitTokenHTB"
    "B22" = "NitroEx"
    "C22" = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
    "E22" = "21NitroExNTX"
    "B23" = "UpBots"
    "C23" = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
    "E23" = "22UpBotsUBXT"
    "B24" = "LEO"
    "C24" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "E24" = "23LEOLEO"
    "B25" = "BTSEToken"
    "C25" = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
    "E25" = "24BTSETokenBTSE"
    "B26" = "BitpandaEcosystemToken"
    "C26" = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
    "E26" = "25BitpandaEcosystemTokenBEST"
    "B27" = "ProBitToken"
    "C27" = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
    "E27" = "26ProBitTokenPROBBestin24h"
    "B41" = "BKEXToken"
    "C41" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "E41" = "40BKEXTokenBKK"
    "B43" = "KickToken"
    "C43" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "E43" = "42KickTokenKICKWorstin24h"
    "E48" = "47BOLOBOLO"
}

foreach ($ref in $textCells.Keys) {
    $ws.Range($ref).Value = $textCells[$ref]
}
